$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 542.94446
$ws.Range("I19").Value = 476.6154
$ws.Range("K19").Value = 476.6154
$ws.Range("M19").Value = -301.6154

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2300.1177
$ws.Range("I40").Value = 1987.5
$ws.Range("J40").Value = 2578
$ws.Range("K40").Value = 1987.5
$ws.Range("L40").Value = 2578
$ws.Range("M40").Value = -1812.5
$ws.Range("N40").Value = -2928

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3495.8333
$ws.Range("J86").Value = 3999.5
$ws.Range("L86").Value = 3999.5
$ws.Range("N86").Value = -6245.5

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4363.273
$ws.Range("I88").Value = 3998.5
$ws.Range("J88").Value = 4444.3335
$ws.Range("K88").Value = 3998.5
$ws.Range("L88").Value = 4444.3335
$ws.Range("M88").Value = -3592.5
$ws.Range("N88").Value = -5256.3335

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3495.8333
$ws.Range("J89").Value = 3999.5
$ws.Range("L89").Value = 19997.5
$ws.Range("N89").Value = -31229.5

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4363.273
$ws.Range("I91").Value = 3998.5
$ws.Range("J91").Value = 4444.3335
$ws.Range("K91").Value = 3998.5
$ws.Range("L91").Value = 4444.3335
$ws.Range("M91").Value = -2594.5
$ws.Range("N91").Value = -7252.3335

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 34811.688
$ws.Range("I106").Value = 35692.08
$ws.Range("J106").Value = 30996.666
$ws.Range("K106").Value = 35692.08
$ws.Range("L106").Value = 30996.666
$ws.Range("M106").Value = -35061.08
$ws.Range("N106").Value = -32258.666

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 37237.75
$ws.Range("I45").Value = 781
$ws.Range("K45").Value = 781
$ws.Range("M45").Value = -404

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1831.7273
$ws.Range("I61").Value = 1831.7273
$ws.Range("K61").Value = 1831.7273
$ws.Range("M61").Value = -1619.7273

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2931.5715
$ws.Range("I102").Value = 1526.6666
$ws.Range("K102").Value = 1526.6666
$ws.Range("M102").Value = 95.33339999999998

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 12568.286
$ws.Range("I110").Value = 13997.6
$ws.Range("J110").Value = 8995
$ws.Range("K110").Value = 13997.6
$ws.Range("L110").Value = 8995
$ws.Range("M110").Value = -11952.6
$ws.Range("N110").Value = -13085

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3781.2
$ws.Range("I122").Value = 3565.7273
$ws.Range("J122").Value = 4373.75
$ws.Range("K122").Value = 10697.1819
$ws.Range("L122").Value = 13121.25
$ws.Range("M122").Value = -8247.1819
$ws.Range("N122").Value = -18021.25

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1831.7273
$ws.Range("I136").Value = 1831.7273
$ws.Range("K136").Value = 5495.1819
$ws.Range("M136").Value = -2945.1819

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5320.8335
$ws.Range("I20").Value = 7376.5
$ws.Range("K20").Value = 7376.5
$ws.Range("M20").Value = -7129.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1430.1177
$ws.Range("I134").Value = 1457
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4371
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1836
$ws.Range("N134").Value = -8070

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4620.2354
$ws.Range("I31").Value = 4102.4287
$ws.Range("K31").Value = 4102.4287
$ws.Range("M31").Value = -3807.4287

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4620.2354
$ws.Range("I34").Value = 4102.4287
$ws.Range("K34").Value = 4102.4287
$ws.Range("M34").Value = -3900.4287

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2051.7693
$ws.Range("I132").Value = 1797.5454
$ws.Range("K132").Value = 5392.6362
$ws.Range("M132").Value = -2862.6362

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3473.077
$ws.Range("I134").Value = 3377.4546
$ws.Range("K134").Value = 10132.3638
$ws.Range("M134").Value = -7597.363799999999

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5622.8667
$ws.Range("I70").Value = 4482.6665
$ws.Range("K70").Value = 4482.6665
$ws.Range("M70").Value = -4212.6665

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5622.8667
$ws.Range("I73").Value = 4482.6665
$ws.Range("K73").Value = 4482.6665
$ws.Range("M73").Value = -3546.6665

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3864
$ws.Range("N40").Value = -5272

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4414.467
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812

# LTW row 64
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 21625
$ws.Range("J64").Value = 21625
$ws.Range("L64").Value = 21625
$ws.Range("N64").Value = -22075

# LTW row 67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 21625
$ws.Range("J67").Value = 21625
$ws.Range("L67").Value = 21625
$ws.Range("N67").Value = -23185

# LTW row 70
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 55000
$ws.Range("I70").Value = 55000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 55000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -54730
$ws.Range("N70").ClearContents()

# LTW row 73
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 55000
$ws.Range("I73").Value = 55000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 55000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -54064
$ws.Range("N73").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6017.1665
$ws.Range("I62").Value = 5500
$ws.Range("J62").Value = 6189.5557
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 6189.5557
$ws.Range("M62").Value = -4876
$ws.Range("N62").Value = -7437.5557

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6017.1665
$ws.Range("I65").Value = 5500
$ws.Range("J65").Value = 6189.5557
$ws.Range("K65").Value = 27500
$ws.Range("L65").Value = 30947.7785
$ws.Range("M65").Value = -24380
$ws.Range("N65").Value = -37187.7785

# WVR row 128
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 29999
$ws.Range("J128").Value = 29999
$ws.Range("L128").Value = 29999
$ws.Range("N128").Value = -39959

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1523.0526
$ws.Range("I136").Value = 1552.1666
$ws.Range("J136").Value = 999
$ws.Range("K136").Value = 4656.4998
$ws.Range("L136").Value = 2997
$ws.Range("M136").Value = -2106.4998
$ws.Range("N136").Value = -8097
